$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15:84 shift down to 16:85.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly price-report record.
$ws.Range("A15").Value = 6
$ws.Range("B15").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C15").Value = "Metropolitana"
$ws.Range("D15").Value = 45099
$ws.Range("E15").Value = 13
$ws.Range("F15").Value = 100112035
$ws.Range("G15").Value = "Bruselas (repollito)"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 140
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 17000
$ws.Range("M15").Value = 17000
$ws.Range("N15").Value = "$/malla 15 kilos"
$ws.Range("O15").Value = "Provincia de Quillota"
$ws.Range("P15").Value = 1133
$ws.Range("Q15").Value = 15
$ws.Range("R15").Value = "Hortaliza"
